# The "About" sheet is the active sheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a "California" label next to the title (row 1), and a date stamp
# (12/7/2023, Excel serial 45267) formatted as a short date (built-in
# numFmtId 14 -> format code "mm-dd-yy").
$ws.Range("B1").Value = "California"
$ws.Range("C1").Value = 45267
$ws.Range("C1").NumberFormat = "mm-dd-yy"
